# Apply the ontology.xlsx edits described by the commit:
# "new .ttl from Google sheet has been generated"
#
# Summary of changes applied to Sheet1:
#   1. Row 23 (ontolidar:VelocityAzimuthDisplay): move the trailing
#      citation sentence out of the definition (column E) into the
#      editorial note (column M).
#   2. Insert a brand-new concept row for "ontolidar:All-fiberLidar"
#      right before the old row 38 (PhotonicsModule), shifting
#      everything at/after row 38 down by one.
#   3. Insert a brand-new concept row for "ontolidar:AcoustoOpticModulator"
#      right before the (now shifted) OpticalAmplifier row, shifting
#      everything at/after that point down by one more.
#   4. Update the altLabel of "ontolidar:MultiLidar" (now row 89) to add
#      ", multistatic".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 23 - split the citation sentence out of the definition into
#    the editorial note column.
# ---------------------------------------------------------------------
$ws.Range("E23").Value2 = "VAD is a method of analyzing data from a complete conical scan whereby many closely spaced azimuthal points may be sampled by the lidar, and the data are used to estimate the wind speed at each height using a statistical fitting method."
$ws.Range("M23").Value2 = "The VAD method is described in Lhermitte (1966) and Browning and Wexler (1968)."

# ---------------------------------------------------------------------
# 2. Insert new row 38: ontolidar:All-fiberLidar
# ---------------------------------------------------------------------
$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value2 = "ontolidar:All-fiberLidar"
$ws.Range("B38").Value2 = "Beam transmission system"
$ws.Range("E38").Value2 = "All-fiber transmission Lidar system"
$ws.Range("H38").Value2 = "boolean"
$ws.Range("M38").Value2 = "True if lidar is all-fiber"

# ---------------------------------------------------------------------
# 3. Insert new row 44: ontolidar:AcoustoOpticModulator
#    (after the previous insert, OpticalAmplifier now sits at row 44,
#    so inserting at row 44 places the new concept right before it)
# ---------------------------------------------------------------------
$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value2 = "ontolidar:AcoustoOpticModulator"
$ws.Range("B44").Value2 = "Acousto optic modulator"
$ws.Range("E44").Value2 = "Bragg cell diffracting and shifting frequency of the light using sound waves (radio-frequency). Produces a pulsed output beam."
$ws.Range("F44").Value2 = "AOM"
$ws.Range("G44").Value2 = "ontolidar:Signal modulation"
$ws.Range("M44").Value2 = "Only pulsed lidar"

# ---------------------------------------------------------------------
# 4. ontolidar:MultiLidar altLabel update (old row 87 -> new row 89)
# ---------------------------------------------------------------------
$ws.Range("F89").Value2 = "multilidar, multistatic"
